# test(web)/qa/salesforce/residencial: agregar planes de residencial 3p 100 cambiar teléfono
#
# On the "Plans" sheet, change the service-type / megas selections for the
# first two residential plan rows, and remove the (duplicate) plan rows 3/4:
#   - C5 "Tipo de Servicio": Sin_TotalPlay_TV -> Con_TotalPlay_TV
#   - D5 "Megas":            50 -> 100
#   - C6 "Tipo de Servicio": Sin_TotalPlay_TV -> Con_TotalPlay_TV_Y_Video_Soundbox
#   - Row 7 (Plan/Tipo de Servicio/Megas) cleared
#   - Row 8 (Plan/Tipo de Servicio/Megas) cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

$ws.Range("C5").Value = "Con_TotalPlay_TV"
$ws.Range("D5").Value = 100

$ws.Range("C6").Value = "Con_TotalPlay_TV_Y_Video_Soundbox"

$ws.Range("B7:D8").ClearContents()

$ws.Rows.Item(6).RowHeight = 23.85

$ws.Range("C8").Select()
